$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these as text (matching the source report's inline-string cells)
$ws.Range("O3:P23").NumberFormat = "@"

$ws.Range("O3").Value = "1551"
$ws.Range("P3").Value = "792580"
$ws.Range("O4").Value = "6"
$ws.Range("P4").Value = "7324"
$ws.Range("O5").Value = "1551"
$ws.Range("P5").Value = "742981"
$ws.Range("O6").Value = "1571"
$ws.Range("P6").Value = "2643.4k"
$ws.Range("O7").Value = "4"
$ws.Range("P7").Value = "7057"
$ws.Range("O8").Value = "1600"
$ws.Range("P8").Value = "8935.7k"
$ws.Range("O9").Value = "1670"
$ws.Range("P9").Value = "2133.3k"
$ws.Range("O10").Value = "1539"
$ws.Range("P10").Value = "3520.8k"
$ws.Range("O11").Value = "1751"
$ws.Range("P11").Value = "809861"
$ws.Range("O12").Value = "3"
$ws.Range("P12").Value = "2934"
$ws.Range("O13").Value = "1981"
$ws.Range("P13").Value = "6948.7k"
$ws.Range("O14").Value = "2"
$ws.Range("P14").Value = "2268"
$ws.Range("O15").Value = "67"
$ws.Range("P15").Value = "11134"
$ws.Range("O16").Value = "1995"
$ws.Range("P16").Value = "6763.5k"
$ws.Range("O17").Value = "1896"
$ws.Range("P17").Value = "539134"
$ws.Range("O18").Value = "1988"
$ws.Range("P18").Value = "403324"
$ws.Range("O19").Value = "1543"
$ws.Range("P19").Value = "2605.4k"
$ws.Range("O20").Value = "2"
$ws.Range("P20").Value = "1607"
$ws.Range("O21").Value = "5"
$ws.Range("P21").Value = "6276"
$ws.Range("O22").Value = "1587"
$ws.Range("P22").Value = "2289.2k"
$ws.Range("O23").Value = "1959"
$ws.Range("P23").Value = "380513"
